$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(46, 8).Value = 4528.5
$ws.Cells.Item(46, 9).Value = 3800
$ws.Cells.Item(46, 10).Value = 4584.5386
$ws.Cells.Item(46, 11).Value = 11400
$ws.Cells.Item(46, 12).Value = 13753.6158
$ws.Cells.Item(46, 13).Value = -11281
$ws.Cells.Item(46, 14).Value = -13991.6158

$ws.Cells.Item(60, 8).Value = 4528.5
$ws.Cells.Item(60, 9).Value = 3800
$ws.Cells.Item(60, 10).Value = 4584.5386
$ws.Cells.Item(60, 11).Value = 11400
$ws.Cells.Item(60, 12).Value = 13753.6158
$ws.Cells.Item(60, 13).Value = -10916
$ws.Cells.Item(60, 14).Value = -14721.6158

$ws.Cells.Item(74, 8).Value = 5448.3887
$ws.Cells.Item(74, 9).Value = 4337.4
$ws.Cells.Item(74, 11).Value = 4337.4
$ws.Cells.Item(74, 13).Value = -3401.4

$ws.Cells.Item(76, 8).Value = 9999
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).Value = $null

$ws.Cells.Item(77, 8).Value = 5448.3887
$ws.Cells.Item(77, 9).Value = 4337.4
$ws.Cells.Item(77, 11).Value = 21687
$ws.Cells.Item(77, 13).Value = -17007

$ws.Cells.Item(79, 8).Value = 9999
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).Value = $null

$ws.Cells.Item(86, 8).Value = 2874.75
$ws.Cells.Item(86, 9).Value = 2500
$ws.Cells.Item(86, 10).Value = 3249.5
$ws.Cells.Item(86, 11).Value = 2500
$ws.Cells.Item(86, 12).Value = 3249.5
$ws.Cells.Item(86, 13).Value = -1377
$ws.Cells.Item(86, 14).Value = -5495.5

$ws.Cells.Item(89, 8).Value = 2874.75
$ws.Cells.Item(89, 9).Value = 2500
$ws.Cells.Item(89, 10).Value = 3249.5
$ws.Cells.Item(89, 11).Value = 12500
$ws.Cells.Item(89, 12).Value = 16247.5
$ws.Cells.Item(89, 13).Value = -6884
$ws.Cells.Item(89, 14).Value = -27479.5

$ws.Cells.Item(127, 8).Value = 5379
$ws.Cells.Item(127, 9).Value = 1500
$ws.Cells.Item(127, 10).Value = 6348.75
$ws.Cells.Item(127, 11).Value = 4500
$ws.Cells.Item(127, 12).Value = 19046.25
$ws.Cells.Item(127, 13).Value = 460
$ws.Cells.Item(127, 14).Value = -28966.25

$ws.Cells.Item(129, 8).Value = 11575.714
$ws.Cells.Item(129, 9).Value = 1421.75
$ws.Cells.Item(129, 11).Value = 4265.25
$ws.Cells.Item(129, 13).Value = 734.75

$ws.Cells.Item(135, 8).Value = 20834184
$ws.Cells.Item(135, 9).Value = 924.8333
$ws.Cells.Item(135, 11).Value = 8323.4997
$ws.Cells.Item(135, 13).Value = -5788.4997

$ws.Cells.Item(138, 8).Value = 1907.0333
$ws.Cells.Item(138, 9).Value = 1056.1875
$ws.Cells.Item(138, 11).Value = 3168.5625
$ws.Cells.Item(138, 13).Value = 1971.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5362.2925
$ws.Cells.Item(32, 9).Value = 4861.054
$ws.Cells.Item(32, 11).Value = 4861.054
$ws.Cells.Item(32, 13).Value = -4574.054

$ws.Cells.Item(45, 8).Value = 3453.7693
$ws.Cells.Item(45, 9).Value = 2958.5
$ws.Cells.Item(45, 11).Value = 2958.5
$ws.Cells.Item(45, 13).Value = -2581.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 55322.332
$ws.Cells.Item(86, 9).Value = 41200.285
$ws.Cells.Item(86, 10).Value = 104749.5
$ws.Cells.Item(86, 11).Value = 41200.285
$ws.Cells.Item(86, 12).Value = 104749.5
$ws.Cells.Item(86, 13).Value = -40077.285
$ws.Cells.Item(86, 14).Value = -106995.5

$ws.Cells.Item(89, 8).Value = 55322.332
$ws.Cells.Item(89, 9).Value = 41200.285
$ws.Cells.Item(89, 10).Value = 104749.5
$ws.Cells.Item(89, 11).Value = 206001.425
$ws.Cells.Item(89, 12).Value = 523747.5
$ws.Cells.Item(89, 13).Value = -200385.425
$ws.Cells.Item(89, 14).Value = -534979.5

$ws.Cells.Item(99, 8).Value = 6363.75
$ws.Cells.Item(99, 9).Value = 4444
$ws.Cells.Item(99, 10).Value = 7003.6665
$ws.Cells.Item(99, 11).Value = 4444
$ws.Cells.Item(99, 12).Value = 7003.6665
$ws.Cells.Item(99, 13).Value = -2946
$ws.Cells.Item(99, 14).Value = -9999.666499999999

$ws.Cells.Item(105, 8).Value = 11015.322
$ws.Cells.Item(105, 9).Value = 22175.8
$ws.Cells.Item(105, 10).Value = 5700.8096
$ws.Cells.Item(105, 11).Value = 22175.8
$ws.Cells.Item(105, 12).Value = 5700.8096
$ws.Cells.Item(105, 13).Value = -20428.8
$ws.Cells.Item(105, 14).Value = -9194.809600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10420696
$ws.Cells.Item(31, 9).Value = 3045.8
$ws.Cells.Item(31, 10).Value = 27783448
$ws.Cells.Item(31, 11).Value = 3045.8
$ws.Cells.Item(31, 12).Value = 27783448
$ws.Cells.Item(31, 13).Value = -2750.8
$ws.Cells.Item(31, 14).Value = -27784038

$ws.Cells.Item(34, 8).Value = 10420696
$ws.Cells.Item(34, 9).Value = 3045.8
$ws.Cells.Item(34, 10).Value = 27783448
$ws.Cells.Item(34, 11).Value = 3045.8
$ws.Cells.Item(34, 12).Value = 27783448
$ws.Cells.Item(34, 13).Value = -2843.8
$ws.Cells.Item(34, 14).Value = -27783852

$ws.Cells.Item(99, 8).Value = 4471.48
$ws.Cells.Item(99, 9).Value = 3840.35
$ws.Cells.Item(99, 11).Value = 3840.35
$ws.Cells.Item(99, 13).Value = -2342.35

$ws.Cells.Item(107, 8).Value = 366.52942
$ws.Cells.Item(107, 9).Value = 370.6875
$ws.Cells.Item(107, 11).Value = 370.6875
$ws.Cells.Item(107, 13).Value = 1549.3125

$ws.Cells.Item(126, 8).Value = 4471.48
$ws.Cells.Item(126, 9).Value = 3840.35
$ws.Cells.Item(126, 11).Value = 11521.05
$ws.Cells.Item(126, 13).Value = -9051.049999999999

$ws.Cells.Item(132, 8).Value = 55348.895
$ws.Cells.Item(132, 9).Value = 58157.5
$ws.Cells.Item(132, 11).Value = 174472.5
$ws.Cells.Item(132, 13).Value = -171942.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 61.333332
$ws.Cells.Item(7, 9).Value = 71.59999999999999
$ws.Cells.Item(7, 11).Value = 214.8
$ws.Cells.Item(7, 13).Value = -102.8

$ws.Cells.Item(18, 8).Value = 3990
$ws.Cells.Item(18, 9).Value = 1586
$ws.Cells.Item(18, 11).Value = 4758
$ws.Cells.Item(18, 13).Value = -4589

$ws.Cells.Item(104, 8).Value = 5368
$ws.Cells.Item(104, 10).Value = 5368
$ws.Cells.Item(104, 12).Value = 16104
$ws.Cells.Item(104, 14).Value = -21346

$ws.Cells.Item(121, 8).Value = 1358.1666
$ws.Cells.Item(121, 10).Value = 883.3333
$ws.Cells.Item(121, 12).Value = 2649.9999
$ws.Cells.Item(121, 14).Value = -5269.9999

$ws.Cells.Item(133, 8).Value = 2772.8
$ws.Cells.Item(133, 9).Value = 1449
$ws.Cells.Item(133, 11).Value = 4347
$ws.Cells.Item(133, 13).Value = 713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1700.6666
$ws.Cells.Item(113, 9).Value = 1700.6666
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1700.6666
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = $null
$ws.Cells.Item(113, 14).Value = 469.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5692.591
$ws.Cells.Item(7, 9).Value = 4737.6665
$ws.Cells.Item(7, 10).Value = 6838.5
$ws.Cells.Item(7, 11).Value = 4737.6665
$ws.Cells.Item(7, 12).Value = 6838.5
$ws.Cells.Item(7, 13).Value = -4625.6665
$ws.Cells.Item(7, 14).Value = -7062.5

$ws.Cells.Item(40, 8).Value = 7854.4287
$ws.Cells.Item(40, 9).Value = 6998
$ws.Cells.Item(40, 10).Value = 9995.5
$ws.Cells.Item(40, 11).Value = 6998
$ws.Cells.Item(40, 12).Value = 9995.5
$ws.Cells.Item(40, 13).Value = -6862
$ws.Cells.Item(40, 14).Value = -10267.5

$ws.Cells.Item(55, 8).Value = 797.25
$ws.Cells.Item(55, 9).Value = 539.3333
$ws.Cells.Item(55, 11).Value = 539.3333
$ws.Cells.Item(55, 13).Value = -366.3333

$ws.Cells.Item(68, 8).Value = 5358.3
$ws.Cells.Item(68, 9).Value = 4198.5
$ws.Cells.Item(68, 11).Value = 4198.5
$ws.Cells.Item(68, 13).Value = -3449.5

$ws.Cells.Item(71, 8).Value = 5358.3
$ws.Cells.Item(71, 9).Value = 4198.5
$ws.Cells.Item(71, 11).Value = 20992.5
$ws.Cells.Item(71, 13).Value = -17248.5

$ws.Cells.Item(82, 8).Value = 2387.7693
$ws.Cells.Item(82, 9).Value = 2489
$ws.Cells.Item(82, 10).Value = 2324.5
$ws.Cells.Item(82, 11).Value = 2489
$ws.Cells.Item(82, 12).Value = 2324.5
$ws.Cells.Item(82, 13).Value = -2128
$ws.Cells.Item(82, 14).Value = -3046.5

$ws.Cells.Item(85, 8).Value = 2387.7693
$ws.Cells.Item(85, 9).Value = 2489
$ws.Cells.Item(85, 10).Value = 2324.5
$ws.Cells.Item(85, 11).Value = 2489
$ws.Cells.Item(85, 12).Value = 2324.5
$ws.Cells.Item(85, 13).Value = -1241
$ws.Cells.Item(85, 14).Value = -4820.5

$ws.Cells.Item(93, 8).Value = 1169289.1
$ws.Cells.Item(93, 9).Value = 2325
$ws.Cells.Item(93, 10).Value = 2169544.2
$ws.Cells.Item(93, 11).Value = 2325
$ws.Cells.Item(93, 12).Value = 2169544.2
$ws.Cells.Item(93, 13).Value = -1077
$ws.Cells.Item(93, 14).Value = -2172040.2

$ws.Cells.Item(122, 8).Value = 2911580.5
$ws.Cells.Item(122, 9).Value = 3999.4285
$ws.Cells.Item(122, 10).Value = 8339065
$ws.Cells.Item(122, 11).Value = 11998.2855
$ws.Cells.Item(122, 12).Value = 25017195
$ws.Cells.Item(122, 13).Value = -9548.2855
$ws.Cells.Item(122, 14).Value = -25022095

$ws.Cells.Item(126, 8).Value = 5692.591
$ws.Cells.Item(126, 9).Value = 4737.6665
$ws.Cells.Item(126, 10).Value = 6838.5
$ws.Cells.Item(126, 11).Value = 14212.9995
$ws.Cells.Item(126, 12).Value = 20515.5
$ws.Cells.Item(126, 13).Value = -11742.9995
$ws.Cells.Item(126, 14).Value = -25455.5

$ws.Cells.Item(132, 8).Value = 4115.619
$ws.Cells.Item(132, 9).Value = 3567.1428
$ws.Cells.Item(132, 10).Value = 5212.5713
$ws.Cells.Item(132, 11).Value = 10701.4284
$ws.Cells.Item(132, 12).Value = 15637.7139
$ws.Cells.Item(132, 13).Value = -8171.428400000001
$ws.Cells.Item(132, 14).Value = -20697.7139

$ws.Cells.Item(136, 8).Value = 4817.091
$ws.Cells.Item(136, 9).Value = 3498.5
$ws.Cells.Item(136, 11).Value = 10495.5
$ws.Cells.Item(136, 13).Value = -7945.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3569.2856
$ws.Cells.Item(81, 10).Value = 4399.8
$ws.Cells.Item(81, 12).Value = 8799.6
$ws.Cells.Item(81, 14).Value = -10921.6

$ws.Cells.Item(84, 8).Value = 3569.2856
$ws.Cells.Item(84, 10).Value = 4399.8
$ws.Cells.Item(84, 12).Value = 43998
$ws.Cells.Item(84, 14).Value = -54606

$ws.Cells.Item(136, 8).Value = 12881.348
$ws.Cells.Item(136, 9).Value = 7627.2856
$ws.Cells.Item(136, 11).Value = 22881.8568
$ws.Cells.Item(136, 13).Value = -20331.8568
